$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 1799.3334
# Row 42: Eye of the Beholder | Hi-Potion of Dexterity
$ws.Range("H42").Value = 1612.5
$ws.Range("I42").Value = 1612.5
$ws.Range("K42").Value = 4837.5
$ws.Range("M42").Value = -4607.5
# Row 75: Tomes Roam on the Range | Dhalmelskin Codex
$ws.Range("H75").Value = 32000
$ws.Range("J75").Value = 32000
$ws.Range("L75").Value = 32000
$ws.Range("N75").Value = -33872
# Row 78: Field Trip to the Unknown (L) | Dhalmelskin Codex
$ws.Range("H78").Value = 32000
$ws.Range("J78").Value = 32000
$ws.Range("L78").Value = 96000
$ws.Range("N78").Value = -105360
# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 1891.9
$ws.Range("I88").Value = 300
$ws.Range("J88").Value = 2068.7778
$ws.Range("K88").Value = 300
$ws.Range("L88").Value = 2068.7778
$ws.Range("M88").Value = 106
$ws.Range("N88").Value = -2880.7778
# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 1891.9
$ws.Range("I91").Value = 300
$ws.Range("J91").Value = 2068.7778
$ws.Range("K91").Value = 300
$ws.Range("L91").Value = 2068.7778
$ws.Range("M91").Value = 1104
$ws.Range("N91").Value = -4876.7778
# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 6559.2
$ws.Range("J112").Value = 6949
$ws.Range("L112").Value = 20847
$ws.Range("N112").Value = -23063

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 41: Skillet Scandal | White Skillet
$ws.Range("H41").Value = 6529.143
$ws.Range("I41").Value = 551.25
$ws.Range("J41").Value = 14499.667
$ws.Range("K41").Value = 551.25
$ws.Range("L41").Value = 14499.667
$ws.Range("M41").Value = -137.25
$ws.Range("N41").Value = -15327.667
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 5159
$ws.Range("I45").Value = 1948.75
$ws.Range("K45").Value = 1948.75
$ws.Range("M45").Value = -1571.75
# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Range("H63").Value = 3498.25
$ws.Range("I63").Value = 3498.25
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3498.25
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -2812.25
# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Range("H66").Value = 3498.25
$ws.Range("I66").Value = 3498.25
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 17491.25
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -14059.25
# Row 96: The Gauntlet Is Cast | High Steel Gauntlets of Fending
$ws.Range("H96").Value = 38274
$ws.Range("J96").Value = 38274
$ws.Range("L96").Value = 38274
$ws.Range("N96").Value = -43766
# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1188.8
$ws.Range("I110").Value = 986
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 986
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1059
$ws.Range("N110").Value = -6090

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run | Iron Rivets
$ws.Range("H22").Value = 575
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1737.2858
$ws.Range("I134").Value = 1944.3636
$ws.Range("K134").Value = 5833.0908
$ws.Range("M134").Value = -3298.0908

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 38: Knock on Wood | Walnut Macuahuitl
$ws.Range("H38").Value = 18814.4
$ws.Range("J38").Value = 41998
$ws.Range("L38").Value = 41998
$ws.Range("N38").Value = -42752
# Row 46: Flintstone Fight | Walnut Macuahuitl
$ws.Range("H46").Value = 18814.4
$ws.Range("J46").Value = 41998
$ws.Range("L46").Value = 41998
$ws.Range("N46").Value = -42420
# Row 88: Hold on Adamantite | Adamantite Spear
$ws.Range("H88").Value = 38454.375
$ws.Range("J88").Value = 38454.375
$ws.Range("L88").Value = 38454.375
$ws.Range("N88").Value = -39266.375
# Row 91: Spears for Stone Vigilantes (L) | Adamantite Spear
$ws.Range("H91").Value = 38454.375
$ws.Range("J91").Value = 38454.375
$ws.Range("L91").Value = 38454.375
$ws.Range("N91").Value = -41262.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 51: The Perks of Life at Sea | Jerked Beef
$ws.Range("H51").Value = 2600
$ws.Range("I51").Value = 1900
$ws.Range("K51").Value = 5700
$ws.Range("M51").Value = -5240
# Row 92: Oh No Udon | Gyr Abanian Flour
$ws.Range("H92").Value = 464
$ws.Range("I92").Value = 480
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = 1440
$ws.Range("L92").Value = 1200
$ws.Range("M92").Value = -192
$ws.Range("N92").Value = -3696

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 58: The Big Red | Red Coral Necklace
$ws.Range("H58").Value = 50000
$ws.Range("I58").Value = 50000
$ws.Range("K58").Value = 50000
$ws.Range("M58").Value = -49723
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
# Row 123: Workplace Workout | Ametrine Ring of Fending
$ws.Range("H123").Value = 70000
$ws.Range("J123").Value = 70000
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -74900
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 2199.6
$ws.Range("I126").Value = 2666
$ws.Range("K126").Value = 7998
$ws.Range("M126").Value = -5528

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 31: Open to Attack | Goatskin Jacket
$ws.Range("H31").Value = 3280
$ws.Range("J31").Value = 12000
$ws.Range("L31").Value = 12000
$ws.Range("N31").Value = -12496
# Row 32: Men Who Scare Up Goats | Goatskin Targe
$ws.Range("H32").Value = 9771
$ws.Range("I32").Value = 2156.5
$ws.Range("K32").Value = 2156.5
$ws.Range("M32").Value = -1839.5
# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 1081901
$ws.Range("J40").Value = 2505002.5
$ws.Range("L40").Value = 2505002.5
$ws.Range("N40").Value = -2505274.5
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 5499
$ws.Range("I46").Value = 5499
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 5499
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5311
# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 504.3
$ws.Range("I55").Value = 505.375
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 505.375
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -332.375
$ws.Range("N55").Value = -846
# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 3416.1428
$ws.Range("I100").Value = 2882.6
$ws.Range("K100").Value = 2882.6
$ws.Range("M100").Value = -2341.6
# Row 106: If the Shoe Fits | Gazelleskin Boots of Casting
$ws.Range("H106").Value = 10185
$ws.Range("J106").Value = 10185
$ws.Range("L106").Value = 10185
$ws.Range("N106").Value = -12709
# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 564000.5600000001
$ws.Range("I136").Value = 633875.6
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 1901626.8
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1899076.8
$ws.Range("N136").Value = -20100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1802.6552
$ws.Range("I132").Value = 1802.6552
$ws.Range("K132").Value = 5407.9656
$ws.Range("M132").Value = -2877.9656
